# edit.ps1 — reproduce the commit:
#   1. Re-point the three tables (slides 14, 15, 16) from the custom
#      "Table_0" style {6BFB26D1-3B87-4655-B18F-58B556C32B21} to the
#      built-in style {B4A32829-0279-49D0-9095-73992229473D}.
#   2. Swap the colour palette carried by the deck's (only reachable)
#      theme part from "Red Violet" (Integral) to the stock "Office"
#      palette, matching the target ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$newTableStyle = "{B4A32829-0279-49D0-9095-73992229473D}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Theme colours --------------------------------------------------
# ThemeColorScheme index -> scheme slot:
#  1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
function Set-ThemeRgb($colorScheme, [int]$index, [string]$rrggbb) {
    $r = [Convert]::ToInt32($rrggbb.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($rrggbb.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($rrggbb.Substring(4, 2), 16)
    $colorScheme.Colors($index).RGB = ($b * 65536) + ($g * 256) + $r
}

$officePalette = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme
foreach ($index in $officePalette.Keys) {
    Set-ThemeRgb $colorScheme $index $officePalette[$index]
}
